# Delete row 127 ("子供時代にこうやって坂道で遊んだ人はどれくらいいますか？" post),
# which removes that post entirely and shifts all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(127).Delete()
